$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header suffixes to the respective input-file
# format-version suffixes ("_FV2310" / "_FV2404").
$ws.Cells.Replace("_old", "_FV2310") | Out-Null
$ws.Cells.Replace("_new", "_FV2404") | Out-Null

# Turn the header + data range into an Excel Table ("Table1") so the
# renamed headers are also reflected as the table's column names.
$rng = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split below row 1, frozen).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
